# DGSalloc_master.xlsx update:
#  - DGSalloc_key: shift the year-header row from 1987-2012 to 1997-2017
#  - DGSalloc: update P-supplementation figures (shift newer columns left,
#    dropping the oldest two survey columns) and drop the now-empty trailing
#    column
#  - Active sheet moves from "DDGS market share estimates" to "DGSalloc"

$wb = $excel.ActiveWorkbook

$wsKey   = $wb.Worksheets.Item("DGSalloc_key")
$wsAlloc = $wb.Worksheets.Item("DGSalloc")

# ---------------------------------------------------------------------------
# DGSalloc_key: year headers B1:G1 (1987,1992,1997,2002,2007,2012)
# become B1:F1 (1997,2002,2007,2012,2017); column G is cleared so the sheet
# shrinks back down to A1:F20.
# ---------------------------------------------------------------------------
$wsKey.Cells.Item(1,2).Value = 1997
$wsKey.Cells.Item(1,3).Value = 2002
$wsKey.Cells.Item(1,4).Value = 2007
$wsKey.Cells.Item(1,5).Value = 2012
$wsKey.Cells.Item(1,6).Value = 2017
$wsKey.Columns.Item(7).ClearContents()

$wsKey.Range("F3").Select()

# ---------------------------------------------------------------------------
# DGSalloc: columns C:F hold successive survey years; the two oldest
# (previously duplicated) columns are dropped, the remaining values shift
# left, and the new rightmost column reuses the last survey's figure.
# Net effect per row: C <- old E, D <- old F, E <- old F.
# ---------------------------------------------------------------------------
$wsAlloc.Cells.Item(1,3).Value = 33.6
$wsAlloc.Cells.Item(1,4).Value = 33
$wsAlloc.Cells.Item(1,5).Value = 33
$wsAlloc.Cells.Item(2,3).Value = 24.4
$wsAlloc.Cells.Item(2,4).Value = 25.4
$wsAlloc.Cells.Item(2,5).Value = 25.4
$wsAlloc.Cells.Item(3,4).Value = 2.1
$wsAlloc.Cells.Item(3,5).Value = 2.1
$wsAlloc.Cells.Item(4,3).Value = 11.2
$wsAlloc.Cells.Item(4,4).Value = 11.5
$wsAlloc.Cells.Item(4,5).Value = 11.5
$wsAlloc.Cells.Item(5,3).Value = 2.2999999999999998
$wsAlloc.Cells.Item(5,5).Value = 2.4
$wsAlloc.Cells.Item(8,3).Value = 6.8
$wsAlloc.Cells.Item(8,4).Value = 6.8
$wsAlloc.Cells.Item(9,4).Value = 1.5
$wsAlloc.Cells.Item(9,5).Value = 1.5
$wsAlloc.Cells.Item(10,3).Value = 7.4
$wsAlloc.Cells.Item(10,4).Value = 7.3
$wsAlloc.Cells.Item(10,5).Value = 7.3
$wsAlloc.Cells.Item(13,3).Value = 1.4
$wsAlloc.Cells.Item(13,4).Value = 1.1000000000000001
$wsAlloc.Cells.Item(13,5).Value = 1.1000000000000001
$wsAlloc.Cells.Item(14,4).Value = 1.1000000000000001
$wsAlloc.Cells.Item(14,5).Value = 1.1000000000000001
$wsAlloc.Cells.Item(15,3).Value = 7.1
$wsAlloc.Cells.Item(15,4).Value = 7.2
$wsAlloc.Cells.Item(15,5).Value = 7.2

# Column F is now redundant (its values were copied into E above); clear it
# so the used range shrinks from A1:F19 down to A1:E19.
$wsAlloc.Range("F1:F19").ClearContents()

# DGSalloc becomes the active sheet/tab, replacing "DDGS market share
# estimates"; selection lands on H1.
$wsAlloc.Activate()
$wsAlloc.Range("H1").Select()
